$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 30; all existing rows 30-123 shift
# down to 31-124 (dimension grows from A1:R123 to A1:R124).
$ws.Rows(30).Insert()

# Populate the newly inserted row 30 with the new "Poroto granado" record.
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44624
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 100112030
$ws.Range("G30").Value = "Poroto granado"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 22000
$ws.Range("L30").Value = 22000
$ws.Range("M30").Value = 22000
$ws.Range("N30").Value = "`$/saco 25 kilos"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 880
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
